$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix mislabeled "Package" column (0402 Capacitor -> correct package) ---
$ws.Range("F2").Value = "0603 Capacitor"
$ws.Range("F4").Value = "0603 Capacitor"
$ws.Range("F5").Value = "0603 Capacitor"
$ws.Range("F6").Value = "1210 Capacitor"

# --- Fill in missing DIGIKEY part numbers ---
$ws.Range("G5").Value = "490-10477-1-ND"

$ws.Range("G13").Value = "RHM27CECT-ND"
$ws.Range("G14").Value = "RHM49.9CFCT-ND"
$ws.Range("G15").Value = "RHM1.50KHCT-ND"
$ws.Range("G16").Value = "RHM10.0KHCT-ND"
$ws.Range("G17").Value = "RHM15.0KCFCT-ND"
$ws.Range("G18").Value = "RHM10MDCT-ND"

# --- Swap rows 21 and 22 (and carry row 22's distinct formatting along) ---
$ws.Range("A4:G4").Copy()
$ws.Range("A22:G22").PasteSpecial(-4122)

$ws.Range("A21").Value = 6
$ws.Range("B21").Value = "NC7SZ157"
$ws.Range("C21").Value = "NC7SZ157"
$ws.Range("D21").Value = "SC70-6"
$ws.Range("E21").Value = "S7, S8, S9, S10, S11, S12"
$ws.Range("F21").Value = "MUX"
$ws.Range("G21").Value = "NC7SZ157P6XCT-ND"

$ws.Range("A22").Value = 6
$ws.Range("B22").Value = "AP2401MP"
$ws.Range("C22").Value = "AP2401MP"
$ws.Range("D22").Value = "MSOP-8EP"
$ws.Range("E22").Value = "S13, S14, S15, S16, S17, S18"
$ws.Range("F22").Value = "USB Power Switch"
$ws.Range("G22").Value = "AP2401MP-13DICT-ND "

# --- Rotate rows 26/27/28 (26<-28, 27<-26, 28<-27) ---
$ws.Range("A26").Value = 12
$ws.Range("B26").Value = "SN74LVC2T45"
$ws.Range("C26").Value = "SN74LVC2T45"
$ws.Range("D26").Value = "VFSOP-8"
$ws.Range("E26").Value = "U8, U9, U10, U11, U12, U13, U14, U15, U16, U17, U18, U19"
$ws.Range("F26").Value = "Bidirectional Buffer / Level Converter 2 Circuits"
$ws.Range("G26").Value = "296-17014-1-ND"

$ws.Range("A27").Value = 6
$ws.Range("B27").Value = "TS3USB30E"
$ws.Range("C27").Value = "TS3USB30E"
$ws.Range("D27").Value = "UQFN-10"
$ws.Range("E27").Value = "U20, U21, U22, U23, U24, U25"
$ws.Range("F27").Value = "USB Mux with OE control and ESD"
$ws.Range("G27").Value = "296-24684-1-ND"

$ws.Range("A28").Value = 6
$ws.Range("B28").Value = "MCP23008"
$ws.Range("C28").Value = "MCP23008"
$ws.Range("D28").Value = "SOIC-18"
$ws.Range("E28").Value = "U26, U27, U28, U29, U30, U31"
$ws.Range("F28").Value = "GPIO Extender via I2C - 8 Pin"
$ws.Range("G28").Value = "MCP23008T-E/SOCT-ND"

# --- Move the active selection ---
[void]$ws.Range("H39").Select()
